$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to keep its existing text formatting so that
# numeric-looking price strings (e.g. "1.002") are written back as text, exactly
# matching the source data (which stores every Price cell as a string).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.462.88'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '1.747.73'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '324.35'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('D7').Value = '0.4448'
$ws.Range('E7').Value = '  +4.34%  '
$ws.Range('D8').Value = '0.3568'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '0.07480'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '42.01'
$ws.Range('E10').Value = '  -5.10%  '
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('D13').Value = '20.67'
$ws.Range('E13').Value = '  -4.07%  '
$ws.Range('D14').Value = '6.004'
$ws.Range('E14').Value = '  -1.90%  '
$ws.Range('D15').Value = '7.083'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = '1.749.93'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').Value = '92.73'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '0.00001058'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '0.06413'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('D22').Value = '5.799'
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('D23').Value = '27.530.84'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').Value = '11.15'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = '2.095'
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').Value = '162.61'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').Value = '20.40'
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = '1.950.76'
$ws.Range('E28').Value = '  -2.12%  '
$ws.Range('D29').Value = '2.068'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('D30').Value = '125.58'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '1.070'
$ws.Range('E31').Value = '  -8.34%  '
$ws.Range('E32').Value = '  +4.69%  '
$ws.Range('D33').Value = '0.09038'
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '5.488'
$ws.Range('E34').Value = '  -3.40%  '
$ws.Range('D35').Value = '11.90'
$ws.Range('E35').Value = '  -5.85%  '
$ws.Range('D36').Value = '0.02281'
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('D37').Value = '0.2091'
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06003'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.6345'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = '4.924'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('D42').Value = '1.379'
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('D43').Value = '7.723'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').Value = '13.14'
$ws.Range('E44').Value = '  -3.66%  '
$ws.Range('D45').Value = '3.713'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').Value = '0.5878'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('D47').Value = '121.57'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('D48').Value = '1.944'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').Value = '1.141'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('D50').Value = '0.06828'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').Value = '72.12'
$ws.Range('E51').Value = '  -3.15%  '
